$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: change from the "Kandelabersvamp" record to the "Apelticka" record.
$ws.Cells.Item(13, 1).Value2 = 111424396
$ws.Cells.Item(13, 2).Value2 = 89816
$ws.Cells.Item(13, 4).Value2 = "VU"
$ws.Cells.Item(13, 5).Value2 = 1619
$ws.Cells.Item(13, 6).Value2 = "Apelticka"
$ws.Cells.Item(13, 7).Value2 = "Aurantiporus fissilis"
$ws.Cells.Item(13, 8).Value2 = "(Berk. & M.A.Curtis) H.Jahn ex Ryvarden"
$ws.Cells.Item(13, 29).Value2 = "2 fruktkroppar på två träd ett stående dött träd och en låga."
$ws.Cells.Item(13, 36).Value2 = "asp"
$ws.Cells.Item(13, 37).Value2 = "Populus tremula"
$ws.Cells.Item(13, 41).Value2 = "Populus tremula"

# Row 14: change from the "Apelticka" record to the "Kandelabersvamp" record.
$ws.Cells.Item(14, 1).Value2 = 111424406
$ws.Cells.Item(14, 2).Value2 = 90151
$ws.Cells.Item(14, 4).Value2 = "NT"
$ws.Cells.Item(14, 5).Value2 = 366
$ws.Cells.Item(14, 6).Value2 = "Kandelabersvamp"
$ws.Cells.Item(14, 7).Value2 = "Artomyces pyxidatus"
$ws.Cells.Item(14, 8).Value2 = "(Pers.) Jülich"
$ws.Cells.Item(14, 29).ClearContents()
$ws.Cells.Item(14, 36).ClearContents()
$ws.Cells.Item(14, 37).ClearContents()
$ws.Cells.Item(14, 41).ClearContents()
